$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Size/Pieces shift right to C/D)
$ws.Columns("B").Insert()

# Header for the new GRADE column
$ws.Range("B1").Value = "GRADE"

# Grade values for each line item (Bill in A, Grade in new B, Size in C, Pieces in D)
$grades = @("c1","c2","c3","c1","c2","c3","c1","c2","c2","c3","c2","c2","c1","c2","c2","c3","c1","c2","c2","c3")

for ($i = 0; $i -lt $grades.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $grades[$i]
}

$ws.Range("B4:B6").Select()
